# Episode 04 metabolism workbook update
# - replaces the "dry weight" KEGG/metabolite summary with the
#   fresh-weight leaf-rosette dataset (new headers, new KEGG ids,
#   reworked Starch/Sucrose/Chlorophyll measurements).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: column headers (set first to mirror the author's edit order) ---
$ws.Range("H7").Value = "Starch `n(mg/g FW)"
$ws.Range("J7").Value = "Cholorophyll (mg/g FW)"
$ws.Range("I7").Value = "Sucrose (mg/g FW)"

# --- Row 6: section banner + KEGG ids for the two new metabolites ---
$ws.Range("A6").Value = "Metabolites reported per g of fresh weight of 6-week-old plant leaf rosettes"
$ws.Range("I6").Value = "C00089"
$ws.Range("J6").Value = "C01793"

# --- Row 8: sample A1 ---
$ws.Range("G8").Value = 0.1206
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 1.2
$ws.Range("J8").Value = 1.8

# --- Row 9: sample A2 ---
$ws.Range("G9").Value = 0.1275
$ws.Range("H9").Value = 6.5
$ws.Range("I9").Value = 1.1
$ws.Range("J9").Value = 1.6

# --- Row 10: sample A3 ---
$ws.Range("G10").Value = 0.2872
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1.4

# --- Row 11: sample A4 ---
$ws.Range("G11").Value = 0.1524
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 0.6
$ws.Range("J11").Value = 2

# --- Row 12: sample A5 (only biomass known, rest NA) ---
$ws.Range("G12").Value = 0.2035

# --- Rows 13-17: B1-B5 samples, biomass no longer recorded -> NA ---
$ws.Range("G13").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("G15").Value = "NA"
$ws.Range("G16").Value = "NA"
$ws.Range("G17").Value = "NA"

$ws.Range("K11").Select() | Out-Null
